$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating the existing
#    "2022-Q2" sheet (so it inherits identical styles/column layout),
#    then place it right after "总计" and rename it.
# ---------------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item(1)
$oldQ2Sheet = $wb.Worksheets.Item("2022-Q2")
$oldQ2Sheet.Copy($null, $summarySheet)

$newQ3Sheet = $wb.Worksheets.Item(2)
$newQ3Sheet.Name = "2022-Q3"

# Fill in the brand new Q3 fund data
$newQ3Sheet.Range("C2").Value2 = "嘉实全球房地产（QDII）"
$newQ3Sheet.Range("D2").Value = "'0.38"
$newQ3Sheet.Range("E2").Value = "'94.39"
$newQ3Sheet.Range("F2").Value = "'3.59"
$newQ3Sheet.Range("G2").Value = "'0.0136"
$newQ3Sheet.Range("H2").Value2 = 5

# ---------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q3
#    above the existing data, shifting everything else down.
# ---------------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()

# Restore cell formatting for the newly inserted row (copy format from
# the row below, which still holds the old row 2 styling).
$summarySheet.Range("A3").Copy()
$summarySheet.Range("A2").PasteSpecial(-4122)
$summarySheet.Range("B3:D3").Copy()
$summarySheet.Range("B2:D2").PasteSpecial(-4122)

$summarySheet.Range("A2").Value2 = 0
$summarySheet.Range("B2").Value2 = "2022-Q3"
$summarySheet.Range("C2").Value2 = 1
$summarySheet.Range("D2").Value2 = 0.01

# Column A holds a simple sequential row index (0-based) - renumber it
# for every data row now that a new row has been inserted at the top.
$summarySheet.Range("A3").Value2 = 1
$summarySheet.Range("A4").Value2 = 2
$summarySheet.Range("A5").Value2 = 3
$summarySheet.Range("A6").Value2 = 4
